# Apply statement_sub_section / statement_section corrections on the
# "cbs_6" sheet (first sheet, sheet1.xml in the OOXML package).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cbs_6")

# Column F (statement_sub_section): current -> noncurrent
$ws.Range("F8").Value  = "noncurrent"
$ws.Range("F9").Value  = "noncurrent"
$ws.Range("F10").Value = "noncurrent"
$ws.Range("F11").Value = "noncurrent"

# Row 13: statement_section assets -> equity_liabilities,
#         statement_sub_section noncurrent -> current
$ws.Range("E13").Value = "equity_liabilities"
$ws.Range("F13").Value = "current"

# Column F (statement_sub_section): current -> noncurrent
$ws.Range("F20").Value = "noncurrent"
$ws.Range("F21").Value = "noncurrent"

# Column F (statement_sub_section): noncurrent -> equity
$ws.Range("F23").Value = "equity"
$ws.Range("F24").Value = "equity"
